$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Caso1")
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"
$ws.Range("G2").Value = 1.027590274810791
$ws.Range("H2").Value = 1.019942760467529
$ws.Range("I2").Value = 1.020192623774532
$ws.Range("J2").Value = 1.020403351634741
$ws.Range("G3").Value = 1.025810599327087
$ws.Range("H3").Value = 1.021311640739441
$ws.Range("I3").Value = 1.020923241330554
$ws.Range("J3").Value = 1.021254988387227
$ws.Range("G4").Value = 1.025840163230896
$ws.Range("H4").Value = 1.022302865982056
$ws.Range("I4").Value = 1.021639801118448
$ws.Range("J4").Value = 1.022154152393341
$ws.Range("G5").Value = 1.025842189788818
$ws.Range("H5").Value = 1.022504806518555
$ws.Range("I5").Value = 1.021672062462499
$ws.Range("J5").Value = 1.022143915295601
$ws.Range("G6").Value = 1.026240468025208
$ws.Range("H6").Value = 1.023157715797424
$ws.Range("I6").Value = 1.02192421087299
$ws.Range("J6").Value = 1.022460140287876
$ws.Range("G7").Value = 1.026274681091309
$ws.Range("H7").Value = 1.023244023323059
$ws.Range("I7").Value = 1.021933241163622
$ws.Range("J7").Value = 1.022542841732502
$ws.Range("G8").Value = 1.026357293128967
$ws.Range("H8").Value = 1.023131251335144
$ws.Range("I8").Value = 1.021974255805259
$ws.Range("J8").Value = 1.022584833204746
$ws.Range("G9").Value = 1.026337862014771
$ws.Range("H9").Value = 1.023283123970032
$ws.Range("I9").Value = 1.021988076388508
$ws.Range("J9").Value = 1.022589188069105
$ws.Range("G10").Value = 1.02642035484314
$ws.Range("H10").Value = 1.023334503173828
$ws.Range("I10").Value = 1.022066524096482
$ws.Range("J10").Value = 1.022666376084089
$ws.Range("G11").Value = 1.02642834186554
$ws.Range("H11").Value = 1.023093223571777
$ws.Range("I11").Value = 1.02208223828959
$ws.Range("J11").Value = 1.02283438667655
$ws.Range("G12").Value = 1.02669370174408
$ws.Range("H12").Value = 1.023193836212158
$ws.Range("I12").Value = 1.022137369174121
$ws.Range("J12").Value = 1.022767337039113
$ws.Range("G13").Value = 1.02667248249054
$ws.Range("H13").Value = 1.023270845413208
$ws.Range("I13").Value = 1.022149941661986
$ws.Range("J13").Value = 1.022821225225925
$ws.Range("G14").Value = 1.032123923301697
$ws.Range("H14").Value = 1.025158047676086
$ws.Range("I14").Value = 1.023700070941563
$ws.Range("J14").Value = 1.024323923513293
$ws.Range("G15").Value = 1.032360434532166
$ws.Range("H15").Value = 1.025078058242798
$ws.Range("I15").Value = 1.023769287912502
$ws.Range("J15").Value = 1.024388624355197
$ws.Range("G16").Value = 1.025725960731506
$ws.Range("H16").Value = 1.022520303726196
$ws.Range("I16").Value = 1.021640980568508
$ws.Range("J16").Value = 1.02209953777492
$ws.Range("G17").Value = 1.025674939155579
$ws.Range("H17").Value = 1.022416830062866
$ws.Range("I17").Value = 1.021639523861721
$ws.Range("J17").Value = 1.022209756076336
$ws.Range("G18").Value = 1.025647401809692
$ws.Range("H18").Value = 1.022548317909241
$ws.Range("I18").Value = 1.021637508723574
$ws.Range("J18").Value = 1.022130791097879
$ws.Range("G19").Value = 1.025624513626099
$ws.Range("H19").Value = 1.022569298744202
$ws.Range("I19").Value = 1.021622429577075
$ws.Range("J19").Value = 1.022023737430573

$ws = $wb.Worksheets.Item("Caso2")
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"
$ws.Range("G2").Value = 1.029593706130981
$ws.Range("H2").Value = 1.021716117858887
$ws.Range("I2").Value = 1.021859761950581
$ws.Range("J2").Value = 1.022270161658525
$ws.Range("G3").Value = 1.028064727783203
$ws.Range("H3").Value = 1.023112297058105
$ws.Range("I3").Value = 1.022616147704591
$ws.Range("J3").Value = 1.023100392892957
$ws.Range("G4").Value = 1.028113484382629
$ws.Range("H4").Value = 1.024274230003357
$ws.Range("I4").Value = 1.023330375588826
$ws.Range("J4").Value = 1.023942733183503
$ws.Range("G5").Value = 1.028115749359131
$ws.Range("H5").Value = 1.024326920509338
$ws.Range("I5").Value = 1.023364182288714
$ws.Range("J5").Value = 1.023984387516975
$ws.Range("G6").Value = 1.028516530990601
$ws.Range("H6").Value = 1.025005221366882
$ws.Range("I6").Value = 1.023614515168698
$ws.Range("J6").Value = 1.024216698482633
$ws.Range("G7").Value = 1.028550982475281
$ws.Range("H7").Value = 1.024971604347229
$ws.Range("I7").Value = 1.023622743192304
$ws.Range("J7").Value = 1.024308227002621
$ws.Range("G8").Value = 1.028632879257202
$ws.Range("H8").Value = 1.024969816207886
$ws.Range("I8").Value = 1.023662649086124
$ws.Range("J8").Value = 1.024393614381552
$ws.Range("G9").Value = 1.028615355491638
$ws.Range("H9").Value = 1.025191903114319
$ws.Range("I9").Value = 1.023677071303612
$ws.Range("J9").Value = 1.024362998083234
$ws.Range("G10").Value = 1.028699636459351
$ws.Range("H10").Value = 1.025081276893616
$ws.Range("I10").Value = 1.023758743405119
$ws.Range("J10").Value = 1.024445202201605
$ws.Range("G11").Value = 1.02870762348175
$ws.Range("H11").Value = 1.024889588356018
$ws.Range("I11").Value = 1.023773435578957
$ws.Range("J11").Value = 1.024532794952393
$ws.Range("G12").Value = 1.028972625732422
$ws.Range("H12").Value = 1.024939060211182
$ws.Range("I12").Value = 1.0238265034959
$ws.Range("J12").Value = 1.024550132453442
$ws.Range("G13").Value = 1.028951525688171
$ws.Range("H13").Value = 1.025027394294739
$ws.Range("I13").Value = 1.023837885363732
$ws.Range("J13").Value = 1.024615235626698
$ws.Range("G14").Value = 1.034380435943604
$ws.Range("H14").Value = 1.026807904243469
$ws.Range("I14").Value = 1.025371825666343
$ws.Range("J14").Value = 1.025975652039051
$ws.Range("G15").Value = 1.034617304801941
$ws.Range("H15").Value = 1.026871919631958
$ws.Range("I15").Value = 1.025440093339086
$ws.Range("J15").Value = 1.026066292077303
$ws.Range("G16").Value = 1.028000712394714
$ws.Range("H16").Value = 1.024463653564453
$ws.Range("I16").Value = 1.023332421304764
$ws.Range("J16").Value = 1.023930778726935
$ws.Range("G17").Value = 1.027950167655945
$ws.Range("H17").Value = 1.024410486221313
$ws.Range("I17").Value = 1.023331941790739
$ws.Range("J17").Value = 1.023992216214538
$ws.Range("G18").Value = 1.027922868728638
$ws.Range("H18").Value = 1.024537563323975
$ws.Range("I18").Value = 1.023328702466601
$ws.Range("J18").Value = 1.023992385715246
$ws.Range("G19").Value = 1.027900218963623
$ws.Range("H19").Value = 1.024486184120178
$ws.Range("I19").Value = 1.023315533862272
$ws.Range("J19").Value = 1.023907117545605

$ws = $wb.Worksheets.Item("Caso3")
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"
$ws.Range("G2").Value = 1.028146266937256
$ws.Range("H2").Value = 1.020262360572815
$ws.Range("I2").Value = 1.020357067382873
$ws.Range("J2").Value = 1.020799193531275
$ws.Range("G3").Value = 1.026573419570923
$ws.Range("H3").Value = 1.021763682365417
$ws.Range("I3").Value = 1.021121925789466
$ws.Range("J3").Value = 1.021679824218154
$ws.Range("G4").Value = 1.026589870452881
$ws.Range("H4").Value = 1.02286171913147
$ws.Range("I4").Value = 1.021840670341476
$ws.Range("J4").Value = 1.022544465959072
$ws.Range("G5").Value = 1.026592016220093
$ws.Range("H5").Value = 1.022934675216675
$ws.Range("I5").Value = 1.021873419395756
$ws.Range("J5").Value = 1.022502392530441
$ws.Range("G6").Value = 1.026983380317688
$ws.Range("H6").Value = 1.023601055145264
$ws.Range("I6").Value = 1.022126654256739
$ws.Range("J6").Value = 1.02282833494246
$ws.Range("G7").Value = 1.02701723575592
$ws.Range("H7").Value = 1.023615479469299
$ws.Range("I7").Value = 1.022136641153559
$ws.Range("J7").Value = 1.022885486483574
$ws.Range("G8").Value = 1.027099370956421
$ws.Range("H8").Value = 1.02363395690918
$ws.Range("I8").Value = 1.022175798400352
$ws.Range("J8").Value = 1.022984454408288
$ws.Range("G9").Value = 1.02707827091217
$ws.Range("H9").Value = 1.023738980293274
$ws.Range("I9").Value = 1.02218923435227
$ws.Range("J9").Value = 1.022941535338759
$ws.Range("G10").Value = 1.027158260345459
$ws.Range("H10").Value = 1.023633360862732
$ws.Range("I10").Value = 1.022272775235017
$ws.Range("J10").Value = 1.023027962073684
$ws.Range("G11").Value = 1.027166128158569
$ws.Range("H11").Value = 1.023514747619629
$ws.Range("I11").Value = 1.022286631946801
$ws.Range("J11").Value = 1.023161510005593
$ws.Range("G12").Value = 1.02742862701416
$ws.Range("H12").Value = 1.023550868034363
$ws.Range("I12").Value = 1.022340717738838
$ws.Range("J12").Value = 1.023122208192945
$ws.Range("G13").Value = 1.02740752696991
$ws.Range("H13").Value = 1.023656129837036
$ws.Range("I13").Value = 1.022351701328981
$ws.Range("J13").Value = 1.023185610771179
$ws.Range("G14").Value = 1.032793641090393
$ws.Range("H14").Value = 1.025444984436035
$ws.Range("I14").Value = 1.023910220429243
$ws.Range("J14").Value = 1.024646738544106
$ws.Range("G15").Value = 1.033028841018677
$ws.Range("H15").Value = 1.025494456291199
$ws.Range("I15").Value = 1.023981344830107
$ws.Range("J15").Value = 1.024717016145587
$ws.Range("G16").Value = 1.026476860046387
$ws.Range("H16").Value = 1.023042798042297
$ws.Range("I16").Value = 1.021843167464445
$ws.Range("J16").Value = 1.022486876696348
$ws.Range("G17").Value = 1.026426315307617
$ws.Range("H17").Value = 1.022890329360962
$ws.Range("I17").Value = 1.02184105456308
$ws.Range("J17").Value = 1.022565955296159
$ws.Range("G18").Value = 1.0263991355896
$ws.Range("H18").Value = 1.022961139678955
$ws.Range("I18").Value = 1.021837542433886
$ws.Range("J18").Value = 1.022543009370565
$ws.Range("G19").Value = 1.026376247406006
$ws.Range("H19").Value = 1.022940039634705
$ws.Range("I19").Value = 1.021827112642556
$ws.Range("J19").Value = 1.022475566715002
